# Homunkulus "database.xlsx" update — append a newly-written folder-path
# entry (3 path lines joined with line breaks) next to the existing
# " - " separator string on row 1, then move the selection down to C19
# as the next empty slot for a future entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "C:\Users\Tim\Documents\.16151814`nC:\Users\Tim\Documents\.16151814\Vid`nC:\Users\Tim\Documents\.16151814\Vid\Straight"

# Excel recalculates the row's display height after the multi-line entry;
# re-running AutoFit on row 1 settles it back onto the sheet's default
# height instead of leaving an explicit custom height behind.
$ws.Rows("1:1").AutoFit() | Out-Null

# Move the active selection to where the next path entry would be typed.
$ws.Range("C19").Select() | Out-Null
